$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 38/39: Stellar and EnergySwap rows swap places (EnergySwap now row 38, Stellar now row 39),
# each carrying refreshed price/volume data.
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D38" "23.65"
$ws.Range("E38").Value = "  +1.98%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D39" "0.117"
$ws.Range("E39").Value = "  -4.54%  "

# Price (D) and Volume(1h) (E) refresh for all other rows
Set-TextValue "D2" "42.122.97"
$ws.Range("E2").Value = "  -8.89%  "

Set-TextValue "D3" "2.489.79"
$ws.Range("E3").Value = "  -4.31%  "

Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.22%  "

Set-TextValue "D5" "293.50"
$ws.Range("E5").Value = "  -4.22%  "

Set-TextValue "D6" "91.93"
$ws.Range("E6").Value = "  -7.38%  "

Set-TextValue "D7" "0.565"
$ws.Range("E7").Value = "  -5.80%  "

Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  +0.19%  "

Set-TextValue "D9" "0.539"
$ws.Range("E9").Value = "  -6.37%  "

Set-TextValue "D10" "35.62"
$ws.Range("E10").Value = "  -9.06%  "

Set-TextValue "D11" "0.0790"
$ws.Range("E11").Value = "  -5.92%  "

Set-TextValue "D12" "7.54"
$ws.Range("E12").Value = "  -6.48%  "

$ws.Range("E13").Value = "  +0.10%  "

Set-TextValue "D14" "2.881.59"
$ws.Range("E14").Value = "  -4.05%  "

Set-TextValue "D15" "2.502.13"
$ws.Range("E15").Value = "  -3.92%  "

Set-TextValue "D16" "0.853"
$ws.Range("E16").Value = "  -6.70%  "

Set-TextValue "D17" "13.86"
$ws.Range("E17").Value = "  -6.80%  "

Set-TextValue "D18" "42.187.57"
$ws.Range("E18").Value = "  -8.97%  "

Set-TextValue "D19" "0.0₃0949"
$ws.Range("E19").Value = "  -5.48%  "

Set-TextValue "D20" "6.41"
$ws.Range("E20").Value = "  -3.91%  "

Set-TextValue "D21" "12.13"
$ws.Range("E21").Value = "  -5.78%  "

Set-TextValue "D22" "71.79"
$ws.Range("E22").Value = "  +0.92%  "

Set-TextValue "D23" "254.77"
$ws.Range("E23").Value = "  -6.03%  "

Set-TextValue "D24" "2.85"
$ws.Range("E24").Value = "  -5.57%  "

Set-TextValue "D25" "2.07"
$ws.Range("E25").Value = "  -3.81%  "

Set-TextValue "D26" "28.38"
$ws.Range("E26").Value = "  -2.71%  "

$ws.Range("E27").Value = "  -0.01%  "

Set-TextValue "D28" "2.19"
$ws.Range("E28").Value = "  -0.68%  "

Set-TextValue "D29" "9.78"
$ws.Range("E29").Value = "  -7.11%  "

Set-TextValue "D30" "36.20"
$ws.Range("E30").Value = "  -5.84%  "

Set-TextValue "D31" "5.90"
$ws.Range("E31").Value = "  -6.03%  "

Set-TextValue "D32" "3.41"
$ws.Range("E32").Value = "  -6.19%  "

Set-TextValue "D33" "149.40"
$ws.Range("E33").Value = "  -1.22%  "

Set-TextValue "D34" "2.14"
$ws.Range("E34").Value = "  -3.35%  "

$ws.Range("E35").Value = "  -5.86%  "

Set-TextValue "D36" "0.0784"
$ws.Range("E36").Value = "  -5.70%  "

$ws.Range("E37").Value = "  -7.97%  "

Set-TextValue "D40" "16.21"
$ws.Range("E40").Value = "  +2.84%  "

Set-TextValue "D41" "3.37"
$ws.Range("E41").Value = "  -5.61%  "

$ws.Range("E42").Value = "  -7.44%  "

$ws.Range("E43").Value = "  -6.89%  "

Set-TextValue "D44" "1.991.92"
$ws.Range("E44").Value = "  -5.77%  "

Set-TextValue "D45" "0.998"
$ws.Range("E45").Value = "  -0.08%  "

Set-TextValue "D46" "84.53"
$ws.Range("E46").Value = "  -9.10%  "

Set-TextValue "D47" "1.59"
$ws.Range("E47").Value = "  +5.24%  "

Set-TextValue "D48" "8.72"
$ws.Range("E48").Value = "  -8.22%  "

Set-TextValue "D49" "2.735.04"
$ws.Range("E49").Value = "  -4.34%  "

Set-TextValue "D50" "101.23"
$ws.Range("E50").Value = "  -6.39%  "

$ws.Range("E51").Value = "  -8.45%  "
